$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the 9da7e663...md file (Ready for handoff -> Handed back) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row 3 Status + Latest Handback DateTime ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("H3").Value = "2016-03-17 03:00:17"

# --- de-de sheet: row 3 Status + Latest Handback DateTime ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("H3").Value = "2016-03-17 03:00:31"

Write-Host "Done"
